$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labels = @(
    "Law & Legal",
    "Philo & Relig",
    "Agri. & Vet sciences",
    "Hist. & Archaeology",
    "Env. sciences",
    "Built Env. & Design",
    "Comm. Manage. Tourism",
    "Biological sciences",
    "Technology",
    "Education",
    "Earth sciences",
    "Creat. Arts & Writing",
    "Chemical sciences",
    "Studies Human Society",
    "Lang. Comms. & Culture",
    "Pysch. & Cognit. sciences",
    "Physical sciences",
    "Economics",
    "Med. & Health Sciences",
    "Engineering",
    "Math sciences",
    "Inf. & Comp. sciences"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $labels[$i]
}
